$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new backlog entries (rows 24-25) ---

# Row 24: new "User story" summary row - copy formatting from the previous
# "User story" row (row 22) which carries the darker fill (style index 2).
$ws.Range("A22:C22").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)  # xlPasteFormats

# Row 25: new "Task" detail row - copy formatting from the previous "Task"
# row (row 23) which carries the yellow fill (style index 3).
$ws.Range("A23:C23").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Fill in the new content
$ws.Range("A24").Value = "User story: Preventing Duplicate Joker Calls"
$ws.Range("B24").Value = 3
$ws.Range("A25").Value = "Task: Implement feature to prevent a player from calling Joker if all Joker cards are revealed"
$ws.Range("B25").Value = 3

# --- Widen the first column to fit the longer text ---
$ws.Columns.Item(1).ColumnWidth = 76.17

# --- Update the selected cell shown when the workbook is opened ---
$ws.Range("B14").Select()
